$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("model_8_8_1", -0.2073650657144932, -4.082543489917939, -14.984180585279, -5.019154632058409, 1.336197257041931, 6.221911907196045, 2.0813148021698, 4.273395538330078)
    ,@("model_8_8_2", -0.1698322669258032, -3.895463233426291, -14.72983849203819, -4.826427245002725, 1.294659614562988, 5.992893218994141, 2.048196792602539, 4.136566162109375)
    ,@("model_8_8_13", 0.1025065566990301, -2.460614027406781, -13.70306739175491, -3.428008817104506, 0.9932608008384705, 4.236390113830566, 1.914499998092651, 3.143736124038696)
    ,@("model_8_8_14", 0.1025065653031366, -2.460613978717051, -13.70306739175491, -3.428008772658157, 0.9932608008384705, 4.236390113830566, 1.914499998092651, 3.143736124038696)
    ,@("model_8_8_12", 0.1025065690947424, -2.46061395809154, -13.70306739175491, -3.42800875383019, 0.9932608008384705, 4.236390113830566, 1.914499998092651, 3.143736124038696)
    ,@("model_8_8_22", 0.1025065700689221, -2.460614050695011, -13.70306743913966, -3.428008772658157, 0.9932608008384705, 4.236390113830566, 1.914499998092651, 3.143736124038696)
    ,@("model_8_8_20", 0.1025065700689221, -2.460614050695011, -13.70306743913966, -3.428008772658157, 0.9932608008384705, 4.236390113830566, 1.914499998092651, 3.143736124038696)
    ,@("model_8_8_19", 0.1025065700689221, -2.460614050695011, -13.70306743913966, -3.428008772658157, 0.9932608008384705, 4.236390113830566, 1.914499998092651, 3.143736124038696)
    ,@("model_8_8_18", 0.1025065700689221, -2.460614050695011, -13.70306743913966, -3.428008772658157, 0.9932608008384705, 4.236390113830566, 1.914499998092651, 3.143736124038696)
    ,@("model_8_8_17", 0.1025065700689221, -2.460614050695011, -13.70306743913966, -3.428008772658157, 0.9932608008384705, 4.236390113830566, 1.914499998092651, 3.143736124038696)
    ,@("model_8_8_16", 0.1025065700689221, -2.460614050695011, -13.70306743913966, -3.428008772658157, 0.9932608008384705, 4.236390113830566, 1.914499998092651, 3.143736124038696)
    ,@("model_8_8_15", 0.1025065700689221, -2.460614050695011, -13.70306743913966, -3.428008772658157, 0.9932608008384705, 4.236390113830566, 1.914499998092651, 3.143736124038696)
    ,@("model_8_8_21", 0.1025065700689221, -2.460614050695011, -13.70306743913966, -3.428008772658157, 0.9932608008384705, 4.236390113830566, 1.914499998092651, 3.143736124038696)
    ,@("model_8_8_10", 0.1025075051839143, -2.460613978717051, -13.70301063018659, -3.428003873682826, 0.9932597279548645, 4.236390113830566, 1.914492607116699, 3.14373254776001)
    ,@("model_8_8_11", 0.1025075075628997, -2.46061395809154, -13.70301063018659, -3.428003854854858, 0.9932597279548645, 4.236390113830566, 1.914492607116699, 3.14373254776001)
    ,@("model_8_8_8", 0.1025083877149021, -2.460614043636691, -13.70295901835734, -3.427999429047948, 0.9932587742805481, 4.236390113830566, 1.914485931396484, 3.143729448318481)
    ,@("model_8_8_9", 0.1025083877149021, -2.460614043636691, -13.70295901835734, -3.427999429047948, 0.9932587742805481, 4.236390113830566, 1.914485931396484, 3.143729448318481)
    ,@("model_8_8_7", 0.1025084094421467, -2.460613949638462, -13.70295901835734, -3.427999343241803, 0.9932587742805481, 4.236390113830566, 1.914485931396484, 3.143729448318481)
    ,@("model_8_8_4", 0.1025092487592552, -2.460614157246062, -13.70290774984402, -3.427995107875151, 0.9932578206062317, 4.236390113830566, 1.91447925567627, 3.143726587295532)
    ,@("model_8_8_6", 0.1025092613047586, -2.460613949638462, -13.70290774984402, -3.427994918360858, 0.9932578206062317, 4.236390113830566, 1.91447925567627, 3.143726348876953)
    ,@("model_8_8_5", 0.1025092613047586, -2.460613949638462, -13.70290774984402, -3.427994918360858, 0.9932578206062317, 4.236390113830566, 1.91447925567627, 3.143726348876953)
    ,@("model_8_8_3", 0.1025112726041397, -2.460614005724636, -13.70278690263405, -3.427984627449258, 0.9932555556297302, 4.236389636993408, 1.914463520050049, 3.14371919631958)
    ,@("model_8_8_0", 0.6176956984198088, 0.4800979790590307, -4.651335032910462, 0.03765304184699003, 0.4230981767177582, 0.636449933052063, 0.7358655333518982, 0.6832337379455566)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = [double]$row[$c]
    }
}

$ws.Rows("25:26").Delete()
